$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D9").Value = "Elon Musk가 공학도라고? 학위는 물리학이랑 경제학인데?"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/elon-musk-physics-economics/#utm_source=rss&utm_medium=rss&utm_campaign=elon-musk-physics-economics"

$ws.Range("D39").Value = "Must-read Guide to Hypothesis Tests You Will Never Use<"
$ws.Range("E39").Value = "https://a292run.tistory.com/entry/Must-read-Guide-to-Hypothesis-Tests-You-Will-Never-Use-1"

$ws.Range("D41").Value = "Cloud 환경에서의 효율적인 보안 및 인증 관리"
$ws.Range("E41").Value = "http://cloudinsight.net/cloud/cloud-%ed%99%98%ea%b2%bd%ec%97%90%ec%84%9c%ec%9d%98-%ed%9a%a8%ec%9c%a8%ec%a0%81%ec%9d%b8-%eb%b3%b4%ec%95%88-%eb%b0%8f-%ec%9d%b8%ec%a6%9d-%ea%b4%80%eb%a6%ac/"

$ws.Range("D51").Value = "[flask] 구름 IDE의 항상 켜두기 기능을 이용해서 flask 웹 서버 구동하기"
$ws.Range("E51").Value = "https://bskyvision.com/1157"
